$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook reports localization status for a set of source files.
# A new handoff report was generated: the file b0c1d244-...md now has its own
# recorded "Latest Handoff Datetime" (previously it was missing / duplicated
# from f1f42c94-...md), and the two rows for b0c1d244-...md / f1f42c94-...md
# are re-sorted (b0c1d244 now sorts before f1f42c94). The underlying
# hyperlink relationships (their target URLs) stay anchored to the same
# table row position they were on before the sort, so after the swap the
# row 4 hyperlink (still using the old row-4 relationship) now points at
# what used to be the f1f42c94 target, and row 5 points at what used to be
# the b0c1d244 target - exactly mirroring the source diff.
# ---------------------------------------------------------------------------

# ---------- Sheet 1: Overview ----------
$ws1 = $wb.Worksheets.Item("Overview")

# Original (row-anchored, unchanged) hyperlink target URLs:
$ov_rId4_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/e2e/f1f42c94-9838-4800-a9c8-09e155678299.md"
$ov_rId5_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/e2e/b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md"
$ov_rId2_url = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md"
$ov_rId3_url = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/52bcbd4b-d994-4290-9bba-13a260905a83.md"
$ov_rId6_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/.localization-config"

# Update the data values for rows 4 and 5 (swap the two file entries).
$ws1.Range("A4").Value = "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = "f1f42c94-9838-4800-a9c8-09e155678299.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

# Rebuild all hyperlinks on the sheet with the final display text, keeping
# every target address exactly where it was (tied to row position, not to
# the file it now displays).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $ov_rId2_url, "", "", "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $ov_rId3_url, "", "", "52bcbd4b-d994-4290-9bba-13a260905a83.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), $ov_rId4_url, "", "", "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), $ov_rId5_url, "", "", "f1f42c94-9838-4800-a9c8-09e155678299.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), $ov_rId6_url, "", "", ".localization-config")

# ---------- Sheet 2: zh-cn ----------
$ws2 = $wb.Worksheets.Item("zh-cn")

$zh_rId2_url = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md"
$zh_rId3_url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/243e4f2b3f1d1ba88ad213aee5513dbe58c54a03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.7cfd6499ef4b5ef1dbd128f298bb1e949efd6c25.zh-cn.xlf"
$zh_rId4_url = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/52bcbd4b-d994-4290-9bba-13a260905a83.md"
$zh_rId5_url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/243e4f2b3f1d1ba88ad213aee5513dbe58c54a03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/52bcbd4b-d994-4290-9bba-13a260905a83.5a3aff9d203e6d6cb7e2f7feb8bc6241b16153d9.zh-cn.xlf"
$zh_rId6_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/e2e/f1f42c94-9838-4800-a9c8-09e155678299.md"
$zh_rId7_url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c0ac29f7ce59544455ece23e756d85e4b5b34ae/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f1f42c94-9838-4800-a9c8-09e155678299.5beaf50d6a13e94b3a94eb37a19ec761c68ec3f6.zh-cn.xlf"
$zh_rId8_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/e2e/b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md"
$zh_rId9_url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c0ac29f7ce59544455ece23e756d85e4b5b34ae/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b0c1d244-ec81-4b9e-975c-6d1bf13868a0.dae505c34529038d520eb6547d02af748667bc3a.zh-cn.xlf"
$zh_rId10_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/.localization-config"

$ws2.Range("A4").Value = "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.dae505c34529038d520eb6547d02af748667bc3a.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-03 09:30:03"

$ws2.Range("A5").Value = "f1f42c94-9838-4800-a9c8-09e155678299.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "f1f42c94-9838-4800-a9c8-09e155678299.5beaf50d6a13e94b3a94eb37a19ec761c68ec3f6.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-03-03 09:29:05"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $zh_rId2_url, "", "", "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zh_rId3_url, "", "", "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.7cfd6499ef4b5ef1dbd128f298bb1e949efd6c25.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $zh_rId4_url, "", "", "52bcbd4b-d994-4290-9bba-13a260905a83.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), $zh_rId5_url, "", "", "52bcbd4b-d994-4290-9bba-13a260905a83.5a3aff9d203e6d6cb7e2f7feb8bc6241b16153d9.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), $zh_rId6_url, "", "", "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), $zh_rId7_url, "", "", "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.dae505c34529038d520eb6547d02af748667bc3a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), $zh_rId8_url, "", "", "f1f42c94-9838-4800-a9c8-09e155678299.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), $zh_rId9_url, "", "", "f1f42c94-9838-4800-a9c8-09e155678299.5beaf50d6a13e94b3a94eb37a19ec761c68ec3f6.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), $zh_rId10_url, "", "", ".localization-config")

# ---------- Sheet 3: de-de ----------
$ws3 = $wb.Worksheets.Item("de-de")

$de_rId2_url = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md"
$de_rId3_url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac8e341a0ca65744383c927fbaebfbf02bf6b4a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.7cfd6499ef4b5ef1dbd128f298bb1e949efd6c25.de-de.xlf"
$de_rId4_url = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/52bcbd4b-d994-4290-9bba-13a260905a83.md"
$de_rId5_url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac8e341a0ca65744383c927fbaebfbf02bf6b4a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/52bcbd4b-d994-4290-9bba-13a260905a83.5a3aff9d203e6d6cb7e2f7feb8bc6241b16153d9.de-de.xlf"
$de_rId6_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/e2e/f1f42c94-9838-4800-a9c8-09e155678299.md"
$de_rId7_url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0c82e118d78e43259f69c80317e71f365dbff40/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f1f42c94-9838-4800-a9c8-09e155678299.5beaf50d6a13e94b3a94eb37a19ec761c68ec3f6.de-de.xlf"
$de_rId8_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/e2e/b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md"
$de_rId9_url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0c82e118d78e43259f69c80317e71f365dbff40/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b0c1d244-ec81-4b9e-975c-6d1bf13868a0.dae505c34529038d520eb6547d02af748667bc3a.de-de.xlf"
$de_rId10_url = "https://github.com/OpenLocalizationTest/oltest/blob/e3db6851ba5e7fd4099c063a5fb534a1ba7876af/.localization-config"

$ws3.Range("A4").Value = "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.dae505c34529038d520eb6547d02af748667bc3a.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-03 09:30:19"

$ws3.Range("A5").Value = "f1f42c94-9838-4800-a9c8-09e155678299.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "f1f42c94-9838-4800-a9c8-09e155678299.5beaf50d6a13e94b3a94eb37a19ec761c68ec3f6.de-de.xlf"
$ws3.Range("D5").Value = "2016-03-03 09:29:18"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $de_rId2_url, "", "", "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), $de_rId3_url, "", "", "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.7cfd6499ef4b5ef1dbd128f298bb1e949efd6c25.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $de_rId4_url, "", "", "52bcbd4b-d994-4290-9bba-13a260905a83.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), $de_rId5_url, "", "", "52bcbd4b-d994-4290-9bba-13a260905a83.5a3aff9d203e6d6cb7e2f7feb8bc6241b16153d9.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), $de_rId6_url, "", "", "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), $de_rId7_url, "", "", "b0c1d244-ec81-4b9e-975c-6d1bf13868a0.dae505c34529038d520eb6547d02af748667bc3a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), $de_rId8_url, "", "", "f1f42c94-9838-4800-a9c8-09e155678299.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), $de_rId9_url, "", "", "f1f42c94-9838-4800-a9c8-09e155678299.5beaf50d6a13e94b3a94eb37a19ec761c68ec3f6.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), $de_rId10_url, "", "", ".localization-config")
